$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.04110066666666667
$ws.Range("H2").Value = 0.123302
$ws.Range("I2").Value = 0.02671259512010182
$ws.Range("J2").Value = 0.02671259512010182
$ws.Range("M2").Value = 3.010057666666667
$ws.Range("N2").Value = 9.030173
$ws.Range("O2").Value = 0.01322668995706902
$ws.Range("P2").Value = 0.01322668995706902
$ws.Range("Q2").Value = 0.1237153768051111
$ws.Range("R2").Value = 1.113438391246
$ws.Range("S2").Value = 0.0003533192136023017
$ws.Range("T2").Value = 0.0003533192136023017

# Row 3
$ws.Range("G3").Value = 0.04110066666666667
$ws.Range("H3").Value = 0.123302
$ws.Range("I3").Value = 0.02671259512010182
$ws.Range("J3").Value = 0.02671259512010182
$ws.Range("O3").Value = 0.08188501082233764
$ws.Range("P3").Value = 0.08188501082233764
$ws.Range("Q3").Value = 0.7659085531948888
$ws.Range("R3").Value = 6.893176978753999
$ws.Range("S3").Value = 0.002187361140502261
$ws.Range("T3").Value = 0.002187361140502261

# Row 4
$ws.Range("G4").Value = 0.04110066666666667
$ws.Range("H4").Value = 0.123302
$ws.Range("I4").Value = 0.02671259512010182
$ws.Range("J4").Value = 0.02671259512010182
$ws.Range("M4").Value = 68.92011666666667
$ws.Range("N4").Value = 206.76035
$ws.Range("O4").Value = 0.3028463623969415
$ws.Range("P4").Value = 0.3028463623969415
$ws.Range("Q4").Value = 2.832662741744445
$ws.Range("R4").Value = 25.4939646757
$ws.Range("S4").Value = 0.008089812262305127
$ws.Range("T4").Value = 0.008089812262305127

# Row 5
$ws.Range("G5").Value = 0.04110066666666667
$ws.Range("H5").Value = 0.123302
$ws.Range("I5").Value = 0.02671259512010182
$ws.Range("J5").Value = 0.02671259512010182
$ws.Range("M5").Value = 8.793419333333334
$ws.Range("N5").Value = 26.380258
$ws.Range("O5").Value = 0.03863973520258026
$ws.Range("P5").Value = 0.03863973520258026
$ws.Range("Q5").Value = 0.3614153968795556
$ws.Range("R5").Value = 3.252738571916
$ws.Range("S5").Value = 0.001032167602014472
$ws.Range("T5").Value = 0.001032167602014472

# Row 6
$ws.Range("G6").Value = 0.04110066666666667
$ws.Range("H6").Value = 0.123302
$ws.Range("I6").Value = 0.02671259512010182
$ws.Range("J6").Value = 0.02671259512010182
$ws.Range("M6").Value = 96.58853933333334
$ws.Range("N6").Value = 289.765618
$ws.Range("O6").Value = 0.4244259760587643
$ws.Range("P6").Value = 0.4244259760587643
$ws.Range("Q6").Value = 3.969853358959556
$ws.Range("R6").Value = 35.728680230636
$ws.Range("S6").Value = 0.0113375192569118
$ws.Range("T6").Value = 0.0113375192569118

# Row 7
$ws.Range("G7").Value = 0.04110066666666667
$ws.Range("H7").Value = 0.123302
$ws.Range("I7").Value = 0.02671259512010182
$ws.Range("J7").Value = 0.02671259512010182
$ws.Range("M7").Value = 31.62744833333333
$ws.Range("N7").Value = 94.882345
$ws.Range("O7").Value = 0.1389762255623074
$ws.Range("P7").Value = 0.1389762255623074
$ws.Range("Q7").Value = 1.299909211465555
$ws.Range("R7").Value = 11.69918290319
$ws.Range("S7").Value = 0.003712415644765863
$ws.Range("T7").Value = 0.003712415644765863

# Row 8
$ws.Range("I8").Value = 0.6362555311831452
$ws.Range("J8").Value = 0.636255531183145
$ws.Range("M8").Value = 3.010057666666667
$ws.Range("N8").Value = 9.030173
$ws.Range("O8").Value = 0.01322668995706902
$ws.Range("P8").Value = 0.01322668995706902
$ws.Range("Q8").Value = 2.946722039949778
$ws.Range("R8").Value = 26.520498359548
$ws.Range("S8").Value = 0.008415554644429722
$ws.Range("T8").Value = 0.008415554644429722

# Row 9
$ws.Range("I9").Value = 0.6362555311831452
$ws.Range("J9").Value = 0.636255531183145
$ws.Range("O9").Value = 0.08188501082233764
$ws.Range("P9").Value = 0.08188501082233764
$ws.Range("S9").Value = 0.05209979105670402
$ws.Range("T9").Value = 0.05209979105670402

# Row 10
$ws.Range("I10").Value = 0.6362555311831452
$ws.Range("J10").Value = 0.636255531183145
$ws.Range("M10").Value = 68.92011666666667
$ws.Range("N10").Value = 206.76035
$ws.Range("O10").Value = 0.3028463623969415
$ws.Range("P10").Value = 0.3028463623969415
$ws.Range("Q10").Value = 67.46994551851111
$ws.Range("R10").Value = 607.2295096666001
$ws.Range("S10").Value = 0.1926876731737493
$ws.Range("T10").Value = 0.1926876731737493

# Row 11
$ws.Range("I11").Value = 0.6362555311831452
$ws.Range("J11").Value = 0.636255531183145
$ws.Range("M11").Value = 8.793419333333334
$ws.Range("N11").Value = 26.380258
$ws.Range("O11").Value = 0.03863973520258026
$ws.Range("P11").Value = 0.03863973520258026
$ws.Range("Q11").Value = 8.608394066000889
$ws.Range("R11").Value = 77.47554659400799
$ws.Range("S11").Value = 0.02458474524609378
$ws.Range("T11").Value = 0.02458474524609377

# Row 12
$ws.Range("I12").Value = 0.6362555311831452
$ws.Range("J12").Value = 0.636255531183145
$ws.Range("M12").Value = 96.58853933333334
$ws.Range("N12").Value = 289.765618
$ws.Range("O12").Value = 0.4244259760587643
$ws.Range("P12").Value = 0.4244259760587643
$ws.Range("Q12").Value = 94.5561876810409
$ws.Range("R12").Value = 851.005689129368
$ws.Range("S12").Value = 0.2700433748451939
$ws.Range("T12").Value = 0.2700433748451939

# Row 13
$ws.Range("I13").Value = 0.6362555311831452
$ws.Range("J13").Value = 0.636255531183145
$ws.Range("M13").Value = 31.62744833333333
$ws.Range("N13").Value = 94.882345
$ws.Range("O13").Value = 0.1389762255623074
$ws.Range("P13").Value = 0.1389762255623074
$ws.Range("Q13").Value = 30.96196465046889
$ws.Range("R13").Value = 278.65768185422
$ws.Range("S13").Value = 0.0884243922169745
$ws.Range("T13").Value = 0.08842439221697447

# Row 14
$ws.Range("G14").Value = 0.5185656666666667
$ws.Range("H14").Value = 1.555697
$ws.Range("I14").Value = 0.3370318736967531
$ws.Range("J14").Value = 0.3370318736967531
$ws.Range("M14").Value = 3.010057666666667
$ws.Range("N14").Value = 9.030173
$ws.Range("O14").Value = 0.01322668995706902
$ws.Range("P14").Value = 0.01322668995706902
$ws.Range("Q14").Value = 1.560912560620111
$ws.Range("R14").Value = 14.048213045581
$ws.Range("S14").Value = 0.004457816099036999
$ws.Range("T14").Value = 0.004457816099036999

# Row 15
$ws.Range("G15").Value = 0.5185656666666667
$ws.Range("H15").Value = 1.555697
$ws.Range("I15").Value = 0.3370318736967531
$ws.Range("J15").Value = 0.3370318736967531
$ws.Range("O15").Value = 0.08188501082233764
$ws.Range("P15").Value = 0.08188501082233764
$ws.Range("Q15").Value = 9.663441294379888
$ws.Range("R15").Value = 86.970971649419
$ws.Range("S15").Value = 0.02759785862513136
$ws.Range("T15").Value = 0.02759785862513136

# Row 16
$ws.Range("G16").Value = 0.5185656666666667
$ws.Range("H16").Value = 1.555697
$ws.Range("I16").Value = 0.3370318736967531
$ws.Range("J16").Value = 0.3370318736967531
$ws.Range("M16").Value = 68.92011666666667
$ws.Range("N16").Value = 206.76035
$ws.Range("O16").Value = 0.3028463623969415
$ws.Range("P16").Value = 0.3028463623969415
$ws.Range("Q16").Value = 35.73960624599445
$ws.Range("R16").Value = 321.65645621395
$ws.Range("S16").Value = 0.1020688769608871
$ws.Range("T16").Value = 0.1020688769608871

# Row 17
$ws.Range("G17").Value = 0.5185656666666667
$ws.Range("H17").Value = 1.555697
$ws.Range("I17").Value = 0.3370318736967531
$ws.Range("J17").Value = 0.3370318736967531
$ws.Range("M17").Value = 8.793419333333334
$ws.Range("N17").Value = 26.380258
$ws.Range("O17").Value = 0.03863973520258026
$ws.Range("P17").Value = 0.03863973520258026
$ws.Range("Q17").Value = 4.559965358869556
$ws.Range("R17").Value = 41.03968822982601
$ws.Range("S17").Value = 0.01302282235447201
$ws.Range("T17").Value = 0.01302282235447201

# Row 18
$ws.Range("G18").Value = 0.5185656666666667
$ws.Range("H18").Value = 1.555697
$ws.Range("I18").Value = 0.3370318736967531
$ws.Range("J18").Value = 0.3370318736967531
$ws.Range("M18").Value = 96.58853933333334
$ws.Range("N18").Value = 289.765618
$ws.Range("O18").Value = 0.4244259760587643
$ws.Range("P18").Value = 0.4244259760587643
$ws.Range("Q18").Value = 50.08750029174956
$ws.Range("R18").Value = 450.7875026257461
$ws.Range("S18").Value = 0.1430450819566586
$ws.Range("T18").Value = 0.1430450819566586

# Row 19
$ws.Range("G19").Value = 0.5185656666666667
$ws.Range("H19").Value = 1.555697
$ws.Range("I19").Value = 0.3370318736967531
$ws.Range("J19").Value = 0.3370318736967531
$ws.Range("M19").Value = 31.62744833333333
$ws.Range("N19").Value = 94.882345
$ws.Range("O19").Value = 0.1389762255623074
$ws.Range("P19").Value = 0.1389762255623074
$ws.Range("Q19").Value = 16.40090882994056
$ws.Range("R19").Value = 147.608179469465
$ws.Range("S19").Value = 0.04683941770056706
$ws.Range("T19").Value = 0.04683941770056706
